# Insert one new weekly record at row 165 ("Fruta / hortaliza, semanal"):
# pushes the existing rows 165-202 down to 166-203 and fills the newly
# opened row 165 with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 165:202 down to 166:203, opening up a blank row 165.
$ws.Rows(165).Insert()

$ws.Range("A165").Value = 7
$ws.Range("B165").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C165").Value = "Ñuble"
$ws.Range("D165").Value = 44754
$ws.Range("E165").Value = 16
$ws.Range("F165").Value = "Fruta"
$ws.Range("G165").Value = 100104
$ws.Range("H165").Value = "Frutos de pepita"
$ws.Range("I165").Value = 100104005
$ws.Range("J165").Value = "Pera"
$ws.Range("K165").Value = "Packham's Triumph"
$ws.Range("L165").Value = "Primera"
$ws.Range("M165").Value = 120
$ws.Range("N165").Value = 8000
$ws.Range("O165").Value = 8500
$ws.Range("P165").Value = 8250
$ws.Range("Q165").Value = "$/caja 16 kilos empedrada"
$ws.Range("R165").Value = "Provincia de Curicó"
$ws.Range("S165").Value = 516
$ws.Range("T165").Value = 16
